$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 75 (Caso 7030 - MONTES DE OCA...) - all rows below shift up by one
$ws.Rows.Item(75).Delete()

# Append the two new "Picada" work orders at the bottom (now rows 80 and 81)

# Row 80
$ws.Cells.Item(80, 1).Value = "'-593"
$ws.Cells.Item(80, 1).ClearFormats()
$ws.Cells.Item(80, 2).Value = "'9/10/2025"
$ws.Cells.Item(80, 2).ClearFormats()
$ws.Cells.Item(80, 3).Value = "Husares 2250"
$ws.Cells.Item(80, 4).Value = "'13"
$ws.Cells.Item(80, 4).ClearFormats()
$ws.Cells.Item(80, 5).Value = "'809642190"
$ws.Cells.Item(80, 5).ClearFormats()
$ws.Cells.Item(80, 6).Value = "Optical Power"
$ws.Cells.Item(80, 7).Value = "Pendiente"
$ws.Cells.Item(80, 8).Value = "Picada"
$ws.Cells.Item(80, 9).Value = 1
$ws.Cells.Item(80, 10).Value = "Cambio"
$ws.Cells.Item(80, 11).Value = "Sin equipos"
$ws.Cells.Item(80, 12).Value = "Pasante"
$ws.Cells.Item(80, 13).Value = -58.443269
$ws.Cells.Item(80, 14).Value = -34.552209
$ws.Cells.Item(80, 15).Value = "Saavedra"
$ws.Cells.Item(80, 16).Value = "Capital Norte"

# Row 81
$ws.Cells.Item(81, 1).Value = "'-594"
$ws.Cells.Item(81, 1).ClearFormats()
$ws.Cells.Item(81, 2).Value = "'9/10/2025"
$ws.Cells.Item(81, 2).ClearFormats()
$ws.Cells.Item(81, 3).Value = "Vidal 1861"
$ws.Cells.Item(81, 4).Value = "'13"
$ws.Cells.Item(81, 4).ClearFormats()
$ws.Cells.Item(81, 5).Value = "'809642175"
$ws.Cells.Item(81, 5).ClearFormats()
$ws.Cells.Item(81, 6).Value = "Optical Power"
$ws.Cells.Item(81, 7).Value = "Pendiente"
$ws.Cells.Item(81, 8).Value = "Picada"
$ws.Cells.Item(81, 9).Value = 1
$ws.Cells.Item(81, 10).Value = "Cambio"
$ws.Cells.Item(81, 11).Value = "Sin equipos"
$ws.Cells.Item(81, 12).Value = "Pasante"
$ws.Cells.Item(81, 13).Value = -58.458298
$ws.Cells.Item(81, 14).Value = -34.566511
$ws.Cells.Item(81, 15).Value = "Colegiales"
$ws.Cells.Item(81, 16).Value = "Capital Norte"
